# Auto-generated edit script: updates TPM-derived NATMI LR-pair metrics
# for the Slit2-Robo1 sheet (ECs ligand-expressing-cell count 2 -> 3, with
# all downstream specificity/weight columns recomputed by the source pipeline).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.176022
$ws.Range("H2").Value = 0.5280659999999999
$ws.Range("I2").Value = 0.03293066697281707
$ws.Range("J2").Value = 0.03293066697281707
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3465496666666667
$ws.Range("N2").Value = 1.039649
$ws.Range("O2").Value = 0.008996151488293185
$ws.Range("P2").Value = 0.008996151488293185
$ws.Range("Q2").Value = 0.061000365426
$ws.Range("R2").Value = 0.549003288834
$ws.Range("S2").Value = 0.0002962492686979956
$ws.Range("T2").Value = 0.0002962492686979955
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.176022
$ws.Range("H3").Value = 0.5280659999999999
$ws.Range("I3").Value = 0.03293066697281707
$ws.Range("J3").Value = 0.03293066697281707
$ws.Range("N3").Value = 66.23320799999999
$ws.Range("O3").Value = 0.5731203249593199
$ws.Range("P3").Value = 0.5731203249593199
$ws.Range("Q3").Value = 3.886167246191999
$ws.Range("R3").Value = 34.97550521572799
$ws.Range("S3").Value = 0.01887323455658806
$ws.Range("T3").Value = 0.01887323455658806
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.176022
$ws.Range("H4").Value = 0.5280659999999999
$ws.Range("I4").Value = 0.03293066697281707
$ws.Range("J4").Value = 0.03293066697281707
$ws.Range("M4").Value = 16.08941833333333
$ws.Range("N4").Value = 48.268255
$ws.Range("O4").Value = 0.4176683996767803
$ws.Range("P4").Value = 0.4176683996767803
$ws.Range("Q4").Value = 2.83209159387
$ws.Range("R4").Value = 25.48882434483
$ws.Range("S4").Value = 0.01375409897482551
$ws.Range("T4").Value = 0.01375409897482551
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.176022
$ws.Range("H5").Value = 0.5280659999999999
$ws.Range("I5").Value = 0.03293066697281707
$ws.Range("J5").Value = 0.03293066697281707
$ws.Range("M5").Value = 0.008287000000000001
$ws.Range("N5").Value = 0.024861
$ws.Range("O5").Value = 0.0002151238756065334
$ws.Range("P5").Value = 0.0002151238756065334
$ws.Range("Q5").Value = 0.001458694314
$ws.Range("R5").Value = 0.013128248826
$ws.Range("S5").Value = 0.000007084172705500478
$ws.Range("T5").Value = 0.000007084172705500478
$ws.Range("I6").Value = 0.8002039325901205
$ws.Range("J6").Value = 0.8002039325901203
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.3465496666666667
$ws.Range("N6").Value = 1.039649
$ws.Range("O6").Value = 0.008996151488293185
$ws.Range("P6").Value = 0.008996151488293185
$ws.Range("Q6").Value = 1.482287994458556
$ws.Range("R6").Value = 13.340591950127
$ws.Range("S6").Value = 0.007198755799108671
$ws.Range("T6").Value = 0.00719875579910867
$ws.Range("I7").Value = 0.8002039325901205
$ws.Range("J7").Value = 0.8002039325901203
$ws.Range("N7").Value = 66.23320799999999
$ws.Range("O7").Value = 0.5731203249593199
$ws.Range("P7").Value = 0.5731203249593199
$ws.Range("Q7").Value = 94.43253353090932
$ws.Range("R7").Value = 849.8928017781839
$ws.Range("S7").Value = 0.4586131378797755
$ws.Range("T7").Value = 0.4586131378797755
$ws.Range("I8").Value = 0.8002039325901205
$ws.Range("J8").Value = 0.8002039325901203
$ws.Range("M8").Value = 16.08941833333333
$ws.Range("N8").Value = 48.268255
$ws.Range("O8").Value = 0.4176683996767803
$ws.Range("P8").Value = 0.4176683996767803
$ws.Range("Q8").Value = 68.81885607542945
$ws.Range("R8").Value = 619.369704678865
$ws.Range("S8").Value = 0.3342198959399817
$ws.Range("T8").Value = 0.3342198959399817
$ws.Range("I9").Value = 0.8002039325901205
$ws.Range("J9").Value = 0.8002039325901203
$ws.Range("M9").Value = 0.008287000000000001
$ws.Range("N9").Value = 0.024861
$ws.Range("O9").Value = 0.0002151238756065334
$ws.Range("P9").Value = 0.0002151238756065334
$ws.Range("Q9").Value = 0.03544577240033334
$ws.Range("R9").Value = 0.319011951603
$ws.Range("S9").Value = 0.0001721429712543759
$ws.Range("T9").Value = 0.0001721429712543759
$ws.Range("G10").Value = 0.891934
$ws.Range("H10").Value = 2.675802
$ws.Range("I10").Value = 0.1668654004370625
$ws.Range("J10").Value = 0.1668654004370625
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.3465496666666667
$ws.Range("N10").Value = 1.039649
$ws.Range("O10").Value = 0.008996151488293185
$ws.Range("P10").Value = 0.008996151488293185
$ws.Range("Q10").Value = 0.3090994303886667
$ws.Range("R10").Value = 2.781894873498
$ws.Range("S10").Value = 0.001501146420486519
$ws.Range("T10").Value = 0.001501146420486519
$ws.Range("G11").Value = 0.891934
$ws.Range("H11").Value = 2.675802
$ws.Range("I11").Value = 0.1668654004370625
$ws.Range("J11").Value = 0.1668654004370625
$ws.Range("N11").Value = 66.23320799999999
$ws.Range("O11").Value = 0.5731203249593199
$ws.Range("P11").Value = 0.5731203249593199
$ws.Range("Q11").Value = 19.691883381424
$ws.Range("R11").Value = 177.226950432816
$ws.Range("S11").Value = 0.09563395252295633
$ws.Range("T11").Value = 0.09563395252295633
$ws.Range("G12").Value = 0.891934
$ws.Range("H12").Value = 2.675802
$ws.Range("I12").Value = 0.1668654004370625
$ws.Range("J12").Value = 0.1668654004370625
$ws.Range("M12").Value = 16.08941833333333
$ws.Range("N12").Value = 48.268255
$ws.Range("O12").Value = 0.4176683996767803
$ws.Range("P12").Value = 0.4176683996767803
$ws.Range("Q12").Value = 14.35069925172333
$ws.Range("R12").Value = 129.15629326551
$ws.Range("S12").Value = 0.06969440476197303
$ws.Range("T12").Value = 0.06969440476197304
$ws.Range("G13").Value = 0.891934
$ws.Range("H13").Value = 2.675802
$ws.Range("I13").Value = 0.1668654004370625
$ws.Range("J13").Value = 0.1668654004370625
$ws.Range("M13").Value = 0.008287000000000001
$ws.Range("N13").Value = 0.024861
$ws.Range("O13").Value = 0.0002151238756065334
$ws.Range("P13").Value = 0.0002151238756065334
$ws.Range("Q13").Value = 0.007391457058000001
$ws.Range("R13").Value = 0.06652311352200001
$ws.Range("S13").Value = 0.00003589673164665703
$ws.Range("T13").Value = 0.00003589673164665703
